$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.756.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.522.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("E7").Value = "  -1.72%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.910.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.524.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.842.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  +3.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.42%  "

$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("E35").Value = "  -1.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("E38").Value = "  -3.08%  "

$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.14%  "

$ws.Range("E42").Value = "  +1.11%  "

$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("E45").Value = "  -3.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.020.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.765.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
